$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 53
$ws.Range("F3").Value = 12
$ws.Range("F4").Value = 448
$ws.Range("G4").Value = "不可售"
$ws.Range("F5").Value = 1333
$ws.Range("G5").Value = 68
$ws.Range("F6").Value = 7659
$ws.Range("F9").Value = 2093
$ws.Range("F10").Value = 8455
$ws.Range("G11").Value = 88
$ws.Range("F13").Value = 63
$ws.Range("F14").Value = 5661
$ws.Range("F15").Value = 60
$ws.Range("F16").Value = 2621
$ws.Range("F17").Value = 1140
$ws.Range("F18").Value = 4595
$ws.Range("F24").Value = 3531
$ws.Range("F25").Value = 58
$ws.Range("F27").Value = 22
$ws.Range("F29").Value = 3028
$ws.Range("F30").Value = 38
$ws.Range("F31").Value = 106
$ws.Range("F35").Value = 445
$ws.Range("G38").Value = 55
$ws.Range("F39").Value = 1792
$ws.Range("F42").Value = 20
$ws.Range("F43").Value = 2914

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 126
$ws.Range("F9").Value = 119
$ws.Range("F10").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 267
$ws.Range("F3").Value = 1327

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 267
$ws.Range("F3").Value = 1327
$ws.Range("F4").Value = 12
$ws.Range("F5").Value = 1333
$ws.Range("G5").Value = 68
$ws.Range("F6").Value = 7659
$ws.Range("F9").Value = 2093
$ws.Range("F10").Value = 8455
$ws.Range("G11").Value = 88
$ws.Range("F13").Value = 63
$ws.Range("F14").Value = 5661
$ws.Range("F15").Value = 60
$ws.Range("F16").Value = 2621
$ws.Range("F17").Value = 1140
$ws.Range("F18").Value = 4595
$ws.Range("F22").Value = 126
$ws.Range("F25").Value = 3531
$ws.Range("F26").Value = 58
$ws.Range("F28").Value = 22
$ws.Range("F29").Value = 3028
$ws.Range("F34").Value = 445
$ws.Range("G37").Value = 55
$ws.Range("F39").Value = 1792
$ws.Range("F42").Value = 20
$ws.Range("F43").Value = 2914
$ws.Range("F49").Value = 119
